$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openTickets")

# --- Re-format the three JSON "highlighting" strings in column G (rows 2-4):
# the raw, single-line JSON (with stray zero-width-space separators) is
# replaced by pretty-printed JSON (3-space indent, no space after ':').

$g2 = @"
[
   {
      "start":143,
      "end":147,
      "key":"System"
   },
   {
      "start":104,
      "end":128,
      "key":"Fehlerbeschreibung"
   },
   {
      "start":67,
      "end":77,
      "key":"System"
   }
]
"@

$g3 = @"
[
   {
      "start":229,
      "end":297,
      "key":"Service Anfrage"
   },
   {
      "start":191,
      "end":192,
      "key":"System"
   },
   {
      "start":176,
      "end":191,
      "key":"System"
   },
   {
      "start":129,
      "end":144,
      "key":"System"
   }
]
"@

$g4 = @"
[
   {
      "start":130,
      "end":165,
      "key":"Auslöser"
   },
   {
      "start":37,
      "end":78,
      "key":"Fehlerbeschreibung"
   },
   {
      "start":24,
      "end":36,
      "key":"System"
   }
]
"@

$ws.Range("G2").Value2 = $g2
$ws.Range("G3").Value2 = $g3
$ws.Range("G4").Value2 = $g4

# --- Wrap text on the (now much longer) highlighting cells, matching the
# wrap style already used on column F.
$ws.Range("G2:G4").WrapText = $true

# --- Grow the row heights so the wrapped JSON is fully visible.
$ws.Rows.Item(2).RowHeight = 246.5
$ws.Rows.Item(3).RowHeight = 319
$ws.Rows.Item(4).RowHeight = 246.5

# --- Scroll the sheet view so row 4 is at the top-left of the pane.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 7
